$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 179.47058
$ws.Range("I33").Value = 94.59999999999999
$ws.Range("J33").Value = 300.7143
$ws.Range("K33").Value = 94.59999999999999
$ws.Range("L33").Value = 300.7143
$ws.Range("M33").Value = 134.4
$ws.Range("N33").Value = -758.7143
$ws.Range("H55").Value = 1838976
$ws.Range("I55").Value = 1056.3636
$ws.Range("J55").Value = 5208495.5
$ws.Range("K55").Value = 1056.3636
$ws.Range("L55").Value = 5208495.5
$ws.Range("M55").Value = -842.3635999999999
$ws.Range("N55").Value = -5208923.5
$ws.Range("H70").Value = 1970.4
$ws.Range("I70").Value = 2083.7778
$ws.Range("J70").Value = 950
$ws.Range("K70").Value = 6251.3334
$ws.Range("L70").Value = 2850
$ws.Range("M70").Value = -5981.3334
$ws.Range("N70").Value = -3390
$ws.Range("H73").Value = 1970.4
$ws.Range("I73").Value = 2083.7778
$ws.Range("J73").Value = 950
$ws.Range("K73").Value = 6251.3334
$ws.Range("L73").Value = 2850
$ws.Range("M73").Value = -5315.3334
$ws.Range("N73").Value = -4722
$ws.Range("H98").Value = 44643436
$ws.Range("I98").Value = 44643436
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 44643436
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -44641938
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 44643436
$ws.Range("I122").Value = 44643436
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 133930308
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -133927858
$ws.Range("N122").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 112652
$ws.Range("J45").Value = 1739.6
$ws.Range("L45").Value = 1739.6
$ws.Range("N45").Value = -2493.6
$ws.Range("H74").Value = 12821472
$ws.Range("I74").Value = 16130096
$ws.Range("J74").Value = 555.5
$ws.Range("K74").Value = 16130096
$ws.Range("L74").Value = 555.5
$ws.Range("M74").Value = -16129222
$ws.Range("N74").Value = -2303.5
$ws.Range("H77").Value = 12821472
$ws.Range("I77").Value = 16130096
$ws.Range("J77").Value = 555.5
$ws.Range("K77").Value = 80650480
$ws.Range("L77").Value = 2777.5
$ws.Range("M77").Value = -80646112
$ws.Range("N77").Value = -11513.5
$ws.Range("H104").Value = 27468.75
$ws.Range("J104").Value = 27468.75
$ws.Range("L104").Value = 27468.75
$ws.Range("N104").Value = -34456.75
$ws.Range("H107").Value = 37332.6
$ws.Range("J107").Value = 37332.6
$ws.Range("L107").Value = 37332.6
$ws.Range("N107").Value = -45012.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 29657.715
$ws.Range("J92").Value = 29657.715
$ws.Range("L92").Value = 29657.715
$ws.Range("N92").Value = -34649.715
$ws.Range("H95").Value = 26131
$ws.Range("J95").Value = 26131
$ws.Range("L95").Value = 26131
$ws.Range("N95").Value = -31623
$ws.Range("H99").Value = 1146.1
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 1140.1111
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 1140.1111
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -4136.1111
$ws.Range("H107").Value = 55556010
$ws.Range("I107").Value = 83333700
$ws.Range("J107").Value = 630
$ws.Range("K107").Value = 83333700
$ws.Range("L107").Value = 630
$ws.Range("M107").Value = -83331780
$ws.Range("N107").Value = -4470
$ws.Range("H132").Value = 48686.89
$ws.Range("J132").Value = 48686.89
$ws.Range("L132").Value = 48686.89
$ws.Range("N132").Value = -58806.89

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 39400
$ws.Range("J64").Value = 39400
$ws.Range("L64").Value = 39400
$ws.Range("N64").Value = -39896
$ws.Range("H67").Value = 39400
$ws.Range("J67").Value = 39400
$ws.Range("L67").Value = 39400
$ws.Range("N67").Value = -41116
$ws.Range("H122").Value = 27778692
$ws.Range("I122").Value = 27778692
$ws.Range("K122").Value = 83336076
$ws.Range("M122").Value = -83333626

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3225.37
$ws.Range("J131").Value = 3370.0737
$ws.Range("L131").Value = 10110.2211
$ws.Range("N131").Value = -20190.2211

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1390.3043
$ws.Range("I102").Value = 1287.2354
$ws.Range("K102").Value = 1287.2354
$ws.Range("M102").Value = 334.7646
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 41676256
$ws.Range("I122").Value = 50011330
$ws.Range("K122").Value = 150033990
$ws.Range("M122").Value = -150031540
$ws.Range("N122").ClearContents()
$ws.Range("H134").Value = 12865.2
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 12865.2
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 38595.60000000001
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -43665.60000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1780.9412
$ws.Range("I7").Value = 1734
$ws.Range("K7").Value = 1734
$ws.Range("M7").Value = -1622
$ws.Range("N7").ClearContents()
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42246
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -131232
$ws.Range("H103").Value = 20387.875
$ws.Range("J103").Value = 20387.875
$ws.Range("L103").Value = 20387.875
$ws.Range("N103").Value = -22731.875
$ws.Range("H106").Value = 142874190
$ws.Range("J106").Value = 142874190
$ws.Range("L106").Value = 142874190
$ws.Range("N106").Value = -142876714
$ws.Range("H126").Value = 1780.9412
$ws.Range("I126").Value = 1734
$ws.Range("K126").Value = 5202
$ws.Range("M126").Value = -2732
$ws.Range("N126").ClearContents()
$ws.Range("H135").Value = 30041.334
$ws.Range("J135").Value = 29971.6
$ws.Range("L135").Value = 29971.6
$ws.Range("N135").Value = -40111.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 15000
$ws.Range("I75").Value = 15000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 15000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -14064
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 15000
$ws.Range("I78").Value = 15000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 45000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -40320
$ws.Range("N78").ClearContents()
$ws.Range("H92").Value = 42450
$ws.Range("J92").Value = 42450
$ws.Range("L92").Value = 42450
$ws.Range("N92").Value = -47442
$ws.Range("H122").Value = 20959.5
$ws.Range("I122").Value = 33202.188
$ws.Range("K122").Value = 99606.56400000001
$ws.Range("M122").Value = -97156.56400000001
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 1588.5454
$ws.Range("I126").Value = 841
$ws.Range("J126").Value = 4952.5
$ws.Range("K126").Value = 2523
$ws.Range("L126").Value = 14857.5
$ws.Range("M126").Value = -53
$ws.Range("N126").Value = -19797.5
